$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.107.58"
$ws.Range("E2").Value = "  -3.48%  "
$ws.Range("D3").Value = "3.680.29"
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.94%  "
$ws.Range("D7").Value = "3.680.39"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000239"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.05%  "
$ws.Range("D15").Value = "4.290.11"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "3.671.23"
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").Value = "67.082.25"
$ws.Range("E17").Value = "  -3.62%  "
$ws.Range("E18").Value = "  -4.26%  "
$ws.Range("E19").Value = "  -6.58%  "
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("E22").Value = "  -5.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000140"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.10%  "
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("D34").Value = "3.817.72"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").Value = "3.616.58"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("E36").Value = "  -7.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.986"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.75%  "
$ws.Range("E40").Value = "  -7.58%  "
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "435.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "141.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "39.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.99%  "
$ws.Range("D50").Value = "2.749.39"
$ws.Range("E50").Value = "  -6.59%  "
$ws.Range("E51").Value = "  -5.43%  "
